$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.414.53'
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("D3").Value = '3.145.83'
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.47%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.146.01'
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("E9").Value = '  -3.11%  '
$ws.Range("E10").Value = '  -4.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.35%  '
$ws.Range("E12").Value = '  -4.86%  '
$ws.Range("E13").Value = '  -3.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.36%  '
$ws.Range("D15").Value = '3.625.80'
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").Value = '64.492.99'
$ws.Range("E16").Value = '  -3.07%  '
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '3.141.25'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.715'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.41%  '
$ws.Range("E29").Value = '  -4.01%  '
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("E32").Value = '  -7.02%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '54.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("D39").Value = '0.0₃0762'
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '452.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.15%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0404'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.22%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.125'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '2.903.25'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.276'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.116'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.95%  '
